# Apply the "Changes of 20th June 2022" edit to the RTE Job Creation workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date/time values on row 2.

# P2: RouteWorkStartDate  44687 (2022-05-06) -> 44729 (2022-06-17)
$ws.Range("P2").Value = 44729

# S2: RouteWorkReadyTime  0 (00:00) -> 0.29166666666666669 (07:00)
$ws.Range("S2").Value = 0.29166666666666669

# T2: RouteWorkScheduledEndTime  0.96875 (23:15) -> 0.28125 (06:45)
$ws.Range("T2").Value = 0.28125

# CA2: FirstGenerationDate  44687 (2022-05-06) -> 44729 (2022-06-17)
$ws.Range("CA2").Value = 44729

# CB2: FirstGenerationTime  3.4722222222222224E-2 (00:50) -> 0.29166666666666669 (07:00)
$ws.Range("CB2").Value = 0.29166666666666669

# Update the saved view state (scroll position and active selection).
$ws.Application.ActiveWindow.ScrollColumn = 14
$ws.Range("Q9").Select()
